# ADD results from server
# Update computed result values on the per-year result sheets.

$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("B2").Value = 0.01288510000000001
$ws2025.Range("E2").Value = 0.3709630638625014
$ws2025.Range("I2").Value = 0.3992358833333333
$ws2025.Range("L2").Value = 0.5681807
$ws2025.Range("M2").Value = 0.07758191666666667
$ws2025.Range("N2").Value = 12.7225376573605
$ws2025.Range("O2").Value = 3.412726950301433

$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("B2").Value = 0.06598974870958903
$ws2030.Range("E2").Value = 0.3655698135733221
$ws2030.Range("I2").Value = 0.7285723112601229
$ws2030.Range("L2").Value = 0.2273058887398774
$ws2030.Range("M2").Value = 0.08785800000000005
$ws2030.Range("N2").Value = 8.998982253460911
$ws2030.Range("O2").Value = 6.734518884057666

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Value = 0.1231904873122577
$ws2035.Range("B2").Value = 0.03401982872623448
$ws2035.Range("E2").Value = 0.1592808142308435
$ws2035.Range("I2").Value = 0.4124207654190268
$ws2035.Range("M2").Value = 0.02832824999999994
$ws2035.Range("N2").Value = 7.845576327967116
$ws2035.Range("O2").Value = 0.739739264829506

$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("N2").Value = 0.4591703578189765

$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("A2").Value = 0.1545504724036852
$ws2045.Range("N2").Value = 2.491682046607792
$ws2045.Range("O2").Value = 5.52224645448535
